$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M4").Value = -0.208543246337598
$ws.Range("C14").Value = -0.208543246337598

$ws.Range("T4").Value = 0.166647238906365
$ws.Range("C21").Value = 0.166647238906365

$ws.Range("L5").Value = 0.234363995960803
$ws.Range("D13").Value = 0.234363995960804

$ws.Range("K6").Value = -0.181734497257527
$ws.Range("E12").Value = -0.181734497257527
